# Refactor "logindata" workbook:
#  - rename CaseNo/Case1..Case7 labels to DataNo/Data1..Data7 on the
#    original "logindata" sheet
#  - add a new "Login Test Case" worksheet after it, describing the
#    actual login test cases (steps / expected results)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("logindata")

# --- rename the case labels in column A --------------------------------
$ws1.Range("A1").Value = "DataNo"
$ws1.Range("A2").Value = "Data1"
$ws1.Range("A3").Value = "Data2"
$ws1.Range("A4").Value = "Data3"
$ws1.Range("A5").Value = "Data4"
$ws1.Range("A6").Value = "Data5"
$ws1.Range("A7").Value = "Data6"
$ws1.Range("A8").Value = "Data7"

# selection moved to A8 on the logindata sheet
$ws1.Range("A8").Select()

# --- add the new "Login Test Case" worksheet ----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Login Test Case"

$ws2.Range("A1").Value = "CaseNo"
$ws2.Range("B1").Value = "Steps"
$ws2.Range("C1").Value = "Expected Result"

$ws2.Range("A2").Value = "Case1"
$ws2.Range("B2").Value = "Doğru Kullanıcı adı ve Şifre girilir."
$ws2.Range("C2").Value = "Login olunduğu görülür"

$ws2.Range("A3").Value = "Case2"
$ws2.Range("B3").Value = "Yanlılş Kullanıcı adı ve Şifre girilir."
$ws2.Range("C3").Value = '"Hatalı E-Posta / Şifre. Tekrar Deneyin." mesajı geldiği görülür'

$ws2.Range("A4").Value = "Case3"
$ws2.Range("B4").Value = "Doğru Kullanıcı ad ve Boş şifre girilir."
$ws2.Range("C4").Value = '"Lütfen şifre giriniz." mesajının geldiği görülür'

$ws2.Range("A5").Value = "Case4"
$ws2.Range("B5").Value = "Boş Kullanıcı adı ve Doğru şifre girilir."
$ws2.Range("C5").Value = '"Lütfen email adresinizi giriniz." mesajının geldiği görülür'

# column widths to fit the longer Turkish text (closest the host's
# column-width grid allows to the authored 35.285.../44.426... widths)
$ws2.Columns.Item(2).ColumnWidth = 34.5
$ws2.Columns.Item(3).ColumnWidth = 43.666666666666664

# selection on the new sheet ends at C5
$ws2.Range("C5").Select()

# keep the original sheet active/selected, matching the workbook's
# unchanged bookViews / tabSelected state
$ws1.Activate()
$ws1.Range("A8").Select()
